# Automatic update of files.
# Column C ("Förändrad") holds the date the record was last refreshed.
# Rows 2-11 move from serial 45204 (2023-10-05) to 45207 (2023-10-08).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45204) {
        $cell.Value = 45207
    }
}
